# Auto-generated edit script applying Hades_Profits market-data refresh
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N) per scheduled runner diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 147841
$ws.Range("J6").Value = 221721.5
$ws.Range("L6").Value = 665164.5
$ws.Range("N6").Value = -665388.5

$ws.Range("H112").Value = 2405.2104
$ws.Range("I112").Value = 883.3333
$ws.Range("J112").Value = 2690.5625
$ws.Range("K112").Value = 2649.9999
$ws.Range("L112").Value = 8071.6875
$ws.Range("M112").Value = -1541.9999
$ws.Range("N112").Value = -10287.6875

$ws.Range("H129").Value = 1028.9
$ws.Range("J129").Value = 1043.3158
$ws.Range("L129").Value = 3129.9474
$ws.Range("N129").Value = -13129.9474

$ws.Range("H135").Value = 17421.508
$ws.Range("I135").Value = 21385.854
$ws.Range("J135").Value = 2783.923
$ws.Range("K135").Value = 192472.686
$ws.Range("L135").Value = 25055.307
$ws.Range("M135").Value = -189937.686
$ws.Range("N135").Value = -30125.307

$ws.Range("H138").Value = 3466171.2
$ws.Range("I138").Value = 272112.66
$ws.Range("J138").Value = 11907612
$ws.Range("K138").Value = 816337.98
$ws.Range("L138").Value = 35722836
$ws.Range("M138").Value = -811197.98
$ws.Range("N138").Value = -35733116

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2094.1052
$ws.Range("I2").Value = 2062.3076
$ws.Range("J2").Value = 2163
$ws.Range("K2").Value = 2062.3076
$ws.Range("L2").Value = 2163
$ws.Range("M2").Value = -1949.3076
$ws.Range("N2").Value = -2389

$ws.Range("H32").Value = 14593.959
$ws.Range("I32").Value = 16726.27
$ws.Range("J32").Value = 7733.478
$ws.Range("K32").Value = 16726.27
$ws.Range("L32").Value = 7733.478
$ws.Range("M32").Value = -16439.27
$ws.Range("N32").Value = -8307.477999999999

$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()

$ws.Range("H45").Value = 1750.421
$ws.Range("I45").Value = 1597.2667
$ws.Range("J45").Value = 2324.75
$ws.Range("K45").Value = 1597.2667
$ws.Range("L45").Value = 2324.75
$ws.Range("M45").Value = -1220.2667
$ws.Range("N45").Value = -3078.75

$ws.Range("H61").Value = 100201200
$ws.Range("I61").Value = 143000720
$ws.Range("J61").Value = 335666.66
$ws.Range("K61").Value = 143000720
$ws.Range("L61").Value = 335666.66
$ws.Range("M61").Value = -143000508
$ws.Range("N61").Value = -336090.66

$ws.Range("H110").Value = 2000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 2000
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -6090

$ws.Range("H116").Value = 2094.1052
$ws.Range("I116").Value = 2062.3076
$ws.Range("J116").Value = 2163
$ws.Range("K116").Value = 2062.3076
$ws.Range("L116").Value = 2163
$ws.Range("M116").Value = 231.6923999999999
$ws.Range("N116").Value = -6751

$ws.Range("H132").Value = 84811.24000000001
$ws.Range("I132").Value = 59796.47
$ws.Range("K132").Value = 179389.41
$ws.Range("M132").Value = -176859.41

$ws.Range("H136").Value = 100201200
$ws.Range("I136").Value = 143000720
$ws.Range("J136").Value = 335666.66
$ws.Range("K136").Value = 429002160
$ws.Range("L136").Value = 1006999.98
$ws.Range("M136").Value = -428999610
$ws.Range("N136").Value = -1012099.98

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2094.1052
$ws.Range("I3").Value = 2062.3076
$ws.Range("J3").Value = 2163
$ws.Range("K3").Value = 2062.3076
$ws.Range("L3").Value = 2163
$ws.Range("M3").Value = -1948.3076
$ws.Range("N3").Value = -2391

$ws.Range("H20").Value = 1012.05884
$ws.Range("I20").Value = 710.8889
$ws.Range("J20").Value = 1350.875
$ws.Range("K20").Value = 710.8889
$ws.Range("L20").Value = 1350.875
$ws.Range("M20").Value = -463.8889
$ws.Range("N20").Value = -1844.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 35241.562
$ws.Range("I134").Value = 1823.9524
$ws.Range("J134").Value = 99038.82000000001
$ws.Range("K134").Value = 5471.857199999999
$ws.Range("L134").Value = 297116.46
$ws.Range("M134").Value = -2936.857199999999
$ws.Range("N134").Value = -302186.46

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 17857740
$ws.Range("I6").Value = 28571464
$ws.Range("J6").Value = 1533.3334
$ws.Range("K6").Value = 85714392
$ws.Range("L6").Value = 4600.0002
$ws.Range("M6").Value = -85714279
$ws.Range("N6").Value = -4826.0002

$ws.Range("H131").Value = 12346620
$ws.Range("I131").Value = 71428930
$ws.Range("J131").Value = 1063.4626
$ws.Range("K131").Value = 214286790
$ws.Range("L131").Value = 3190.3878
$ws.Range("M131").Value = -214281750
$ws.Range("N131").Value = -13270.3878

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 25000
$ws.Range("J108").Value = 25000
$ws.Range("L108").Value = 25000
$ws.Range("N108").Value = -32680

$ws.Range("H110").Value = 44801
$ws.Range("J110").Value = 44801
$ws.Range("L110").Value = 44801
$ws.Range("N110").Value = -52981

$ws.Range("H126").Value = 2005.5652
$ws.Range("I126").Value = 1441
$ws.Range("J126").Value = 2306.6667
$ws.Range("K126").Value = 4323
$ws.Range("L126").Value = 6920.000100000001
$ws.Range("M126").Value = -1853
$ws.Range("N126").Value = -11860.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 2852.2
$ws.Range("I35").Value = 2365.25
$ws.Range("J35").Value = 4800
$ws.Range("K35").Value = 2365.25
$ws.Range("L35").Value = 4800
$ws.Range("M35").Value = -2029.25
$ws.Range("N35").Value = -5472

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H132").Value = 35056.848
$ws.Range("I132").Value = 15697.757
$ws.Range("J132").Value = 114644.22
$ws.Range("K132").Value = 47093.271
$ws.Range("L132").Value = 343932.66
$ws.Range("M132").Value = -44563.271
$ws.Range("N132").Value = -348992.66

$ws.Range("H136").Value = 80983.7
$ws.Range("I136").Value = 47273.332
$ws.Range("K136").Value = 141819.996
$ws.Range("M136").Value = -139269.996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 86862.336
$ws.Range("I132").Value = 80419.84
$ws.Range("J132").Value = 101504.37
$ws.Range("K132").Value = 241259.52
$ws.Range("L132").Value = 304513.11
$ws.Range("M132").Value = -238729.52
$ws.Range("N132").Value = -309573.11

$ws.Range("H136").Value = 41435.58
$ws.Range("I136").Value = 34195.168
$ws.Range("J136").Value = 52296.2
$ws.Range("K136").Value = 102585.504
$ws.Range("L136").Value = 156888.6
$ws.Range("M136").Value = -100035.504
$ws.Range("N136").Value = -161988.6

